$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp label
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 15:50"

# Rows 19-64: Canary Islands provinces inserted after Malaga (row 18),
# all subsequent provinces shift down one logical slot, and case totals
# are refreshed with the latest figures for every affected row.
$data = @(
    @(19, "Gran Canaria", 1204, 32, 320, 11),
    @(20, "La Palma", 1204, 32, 47, 2),
    @(21, "Lanzarote", 1204, 32, 42, 3),
    @(22, "Fuerteventura", 1204, 32, 31, 0),
    @(23, "La Gomera", 1204, 32, 7, 0),
    @(24, "El Hierro", 1204, 32, 3, 0),
    @(25, "Asturias", 1158, 78, 1031, 49),
    @(26, "Cantabria", 1100, 24, 1049, 27),
    @(27, "Salamanca", 1078, 181, 787, 110),
    @(28, "Gipuzkoa/Guipuzcoa", 1071, 1626, 684, 41),
    @(29, "Pontevedra", 1060, 153, 1005, 9),
    @(30, "Sevilla", 1052, 18, 1000, 34),
    @(31, "Caceres", 1045, 11, 945, 89),
    @(32, "Granada", 963, 15, 882, 66),
    @(33, "Murcia", 939, 17, 897, 25),
    @(34, "Valladolid", 929, 145, 718, 66),
    @(35, "Leon", 918, 139, 690, 89),
    @(36, "Aragon", 907, 29, 838, 40),
    @(37, "Burgos", 749, 176, 512, 61),
    @(38, "Segovia", 629, 156, 404, 69),
    @(39, "Jaen", 599, 17, 559, 23),
    @(40, "Castello/Castellon", 586, 9, 545, 32),
    @(41, "Guadalajara", 586, 252, 479, 93),
    @(42, "Cordoba", 572, 4, 555, 13),
    @(43, "Soria", 550, 71, 442, 37),
    @(44, "Badajoz", 515, 49, 449, 17),
    @(45, "Cadiz", 507, 10, 484, 13),
    @(46, "Ourense", 458, 153, 415, 8),
    @(47, "Avila", 446, 91, 309, 46),
    @(48, "Palencia", 293, 33, 242, 18),
    @(49, "Lugo", 270, 153, 244, 4),
    @(50, "Cuenca", 268, 252, 187, 62),
    @(51, "Almeria", 223, 6, 203, 14),
    @(52, "Teruel", 222, 14, 196, 12),
    @(53, "Huesca", 215, 19, 185, 11),
    @(54, "Mallorca", 210, 18, 194, 12),
    @(55, "Zamora", 209, 36, 153, 20),
    @(56, "Huelva", 168, 2, 162, 4),
    @(57, "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena", 58, 0, 58, 3),
    @(58, "Melilla", 51, 0, 50, 1),
    @(59, "Ceuta", 29, 0, 28, 1),
    @(60, "Ibiza", 21, 18, 20, 1),
    @(61, "Menorca", 15, 18, 13, 0),
    @(62, "Arroyo de la Luz", 7, 0, 7, 0),
    @(63, "Tenerife", 3, 30, 1056, 36),
    @(64, "Formentera", 0, 10, 0, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
